# Auto-generated edit script: updates numeric values in H:N columns across several
# worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1014.0909
$ws.Range("I2").Value = 145.6875
$ws.Range("J2").Value = 3329.8333
$ws.Range("K2").Value = 145.6875
$ws.Range("L2").Value = 3329.8333
$ws.Range("M2").Value = -32.6875
$ws.Range("N2").Value = -3555.8333
$ws.Range("H4").Value = 229.8
$ws.Range("I4").Value = 185.14285
$ws.Range("J4").Value = 334
$ws.Range("K4").Value = 185.14285
$ws.Range("L4").Value = 334
$ws.Range("M4").Value = -71.14285000000001
$ws.Range("N4").Value = -562
$ws.Range("H9").Value = 2367.4443
$ws.Range("J9").Value = 499.5
$ws.Range("L9").Value = 499.5
$ws.Range("N9").Value = -837.5
$ws.Range("H38").Value = 1367.1578
$ws.Range("I38").Value = 452.45456
$ws.Range("J38").Value = 2624.875
$ws.Range("K38").Value = 1357.36368
$ws.Range("L38").Value = 7874.625
$ws.Range("M38").Value = -985.3636799999999
$ws.Range("N38").Value = -8618.625
$ws.Range("H51").Value = 12075.85
$ws.Range("H55").Value = 217
$ws.Range("I55").Value = 89
$ws.Range("K55").Value = 89
$ws.Range("M55").Value = 125
$ws.Range("H64").Value = 333333340
$ws.Range("I64").Value = 333333340
$ws.Range("K64").Value = 333333340
$ws.Range("M64").Value = -333333092
$ws.Range("H67").Value = 333333340
$ws.Range("I67").Value = 333333340
$ws.Range("K67").Value = 333333340
$ws.Range("M67").Value = -333332482
$ws.Range("H103").Value = 377.14285
$ws.Range("I103").Value = 377.14285
$ws.Range("K103").Value = 1131.42855
$ws.Range("M103").Value = -545.4285500000001
$ws.Range("H137").Value = 12060.363
$ws.Range("I137").Value = 8973.091
$ws.Range("K137").Value = 26919.273
$ws.Range("M137").Value = -24369.273

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 359.125
$ws.Range("J4").Value = 572.25
$ws.Range("L4").Value = 572.25
$ws.Range("N4").Value = -804.25
$ws.Range("H5").Value = 137.27272
$ws.Range("I5").Value = 137.27272
$ws.Range("K5").Value = 137.27272
$ws.Range("M5").Value = -25.27271999999999
$ws.Range("H33").Value = 16851.084
$ws.Range("J33").Value = 3799.3333
$ws.Range("L33").Value = 3799.3333
$ws.Range("N33").Value = -4457.3333
$ws.Range("H34").Value = 33333
$ws.Range("I34").Value = 33333
$ws.Range("K34").Value = 33333
$ws.Range("M34").Value = -33062
$ws.Range("H36").Value = 7150589
$ws.Range("H45").Value = 2833.1667
$ws.Range("I45").Value = 2833.1667
$ws.Range("K45").Value = 2833.1667
$ws.Range("M45").Value = -2456.1667
$ws.Range("H61").Value = 4751.3335
$ws.Range("I61").Value = 5210.4443
$ws.Range("J61").Value = 3833.111
$ws.Range("K61").Value = 5210.4443
$ws.Range("L61").Value = 3833.111
$ws.Range("M61").Value = -4998.4443
$ws.Range("N61").Value = -4257.111
$ws.Range("H63").Value = 353.55554
$ws.Range("I63").Value = 341
$ws.Range("K63").Value = 341
$ws.Range("M63").Value = 345
$ws.Range("H66").Value = 353.55554
$ws.Range("I66").Value = 341
$ws.Range("K66").Value = 1705
$ws.Range("M66").Value = 1727
$ws.Range("H74").Value = 2627.818
$ws.Range("I74").Value = 2766.2
$ws.Range("K74").Value = 2766.2
$ws.Range("M74").Value = -1892.2
$ws.Range("H77").Value = 2627.818
$ws.Range("I77").Value = 2766.2
$ws.Range("K77").Value = 13831
$ws.Range("M77").Value = -9463
$ws.Range("H122").Value = 6500
$ws.Range("I122").Value = 6000
$ws.Range("K122").Value = 18000
$ws.Range("M122").Value = -15550
$ws.Range("H136").Value = 4751.3335
$ws.Range("I136").Value = 5210.4443
$ws.Range("J136").Value = 3833.111
$ws.Range("K136").Value = 15631.3329
$ws.Range("L136").Value = 11499.333
$ws.Range("M136").Value = -13081.3329
$ws.Range("N136").Value = -16599.333

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 137.27272
$ws.Range("I4").Value = 137.27272
$ws.Range("K4").Value = 137.27272
$ws.Range("M4").Value = -22.27271999999999
$ws.Range("H99").Value = 8118.5
$ws.Range("J99").Value = 9053.817999999999
$ws.Range("L99").Value = 9053.817999999999
$ws.Range("N99").Value = -12049.818

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 19235634
$ws.Range("I58").Value = 26319110
$ws.Range("J58").Value = 9054.929
$ws.Range("K58").Value = 26319110
$ws.Range("L58").Value = 9054.929
$ws.Range("M58").Value = -26318907
$ws.Range("N58").Value = -9460.929
$ws.Range("H62").Value = 19998.834
$ws.Range("I62").Value = 22398.6
$ws.Range("K62").Value = 22398.6
$ws.Range("M62").Value = -21774.6
$ws.Range("H65").Value = 19998.834
$ws.Range("I65").Value = 22398.6
$ws.Range("K65").Value = 111993
$ws.Range("M65").Value = -108873
$ws.Range("H132").Value = 5048.2617
$ws.Range("I132").Value = 4582.1333
$ws.Range("K132").Value = 13746.3999
$ws.Range("M132").Value = -11216.3999
$ws.Range("H136").Value = 19235634
$ws.Range("I136").Value = 26319110
$ws.Range("J136").Value = 9054.929
$ws.Range("K136").Value = 78957330
$ws.Range("L136").Value = 27164.787
$ws.Range("M136").Value = -78954780
$ws.Range("N136").Value = -32264.787

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3548580.5
$ws.Range("I4").Value = 3548580.5
$ws.Range("K4").Value = 10645741.5
$ws.Range("M4").Value = -10645629.5
$ws.Range("H60").Value = 900
$ws.Range("J60").Value = 900
$ws.Range("L60").Value = 2700
$ws.Range("N60").Value = -3202
$ws.Range("H68").Value = 557545.9
$ws.Range("I68").Value = 1832.75
$ws.Range("K68").Value = 5498.25
$ws.Range("M68").Value = -4687.25
$ws.Range("H71").Value = 557545.9
$ws.Range("I71").Value = 1832.75
$ws.Range("K71").Value = 16494.75
$ws.Range("M71").Value = -12438.75
$ws.Range("H133").Value = 3460
$ws.Range("I133").Value = 3460
$ws.Range("K133").Value = 10380
$ws.Range("M133").Value = -5320

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 125
$ws.Range("I17").Value = 125
$ws.Range("K17").Value = 125
$ws.Range("M17").Value = 43
$ws.Range("H107").Value = 732.9167
$ws.Range("I107").Value = 663.375
$ws.Range("J107").Value = 872
$ws.Range("K107").Value = 663.375
$ws.Range("L107").Value = 872
$ws.Range("M107").Value = 1256.625
$ws.Range("N107").Value = -4712
$ws.Range("H122").Value = 3597.875
$ws.Range("I122").Value = 2921.375
$ws.Range("K122").Value = 8764.125
$ws.Range("M122").Value = -6314.125
$ws.Range("H132").Value = 15388088
$ws.Range("I132").Value = 24393496
$ws.Range("J132").Value = 3849
$ws.Range("K132").Value = 73180488
$ws.Range("L132").Value = 11547
$ws.Range("M132").Value = -73177958
$ws.Range("N132").Value = -16607

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 4500
$ws.Range("I2").Value = 4500
$ws.Range("K2").Value = 4500
$ws.Range("M2").Value = -4388
$ws.Range("H23").Value = 16683333
$ws.Range("I23").Value = 16683333
$ws.Range("K23").Value = 16683333
$ws.Range("M23").Value = -16683103
$ws.Range("H38").Value = 34665.8
$ws.Range("J38").Value = 35554.332
$ws.Range("L38").Value = 35554.332
$ws.Range("N38").Value = -36374.332
$ws.Range("H46").Value = 41667776
$ws.Range("I46").Value = 1299.5
$ws.Range("J46").Value = 62501012
$ws.Range("K46").Value = 1299.5
$ws.Range("L46").Value = 62501012
$ws.Range("M46").Value = -1111.5
$ws.Range("N46").Value = -62501388
$ws.Range("H50").Value = 31999.666
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 31999.666
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 31999.666
$ws.Range("N50").Value = -33273.666
$ws.Range("M50").ClearContents()
$ws.Range("H68").Value = 1698.9333
$ws.Range("I68").Value = 1498.7273
$ws.Range("J68").Value = 2249.5
$ws.Range("K68").Value = 1498.7273
$ws.Range("L68").Value = 2249.5
$ws.Range("M68").Value = -749.7273
$ws.Range("N68").Value = -3747.5
$ws.Range("H71").Value = 1698.9333
$ws.Range("I71").Value = 1498.7273
$ws.Range("J71").Value = 2249.5
$ws.Range("K71").Value = 7493.636500000001
$ws.Range("L71").Value = 11247.5
$ws.Range("M71").Value = -3749.636500000001
$ws.Range("N71").Value = -18735.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 26273.727
$ws.Range("J14").Value = 25075.666
$ws.Range("L14").Value = 25075.666
$ws.Range("N14").Value = -25411.666
$ws.Range("H34").Value = 32666.334
$ws.Range("J34").Value = 31999.666
$ws.Range("L34").Value = 31999.666
$ws.Range("N34").Value = -32405.666
$ws.Range("H42").Value = 47048.5
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 47048.5
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 47048.5
$ws.Range("N42").Value = -47804.5
$ws.Range("M42").ClearContents()
$ws.Range("H49").Value = 32650.285
$ws.Range("I49").Value = 32249.75
$ws.Range("J49").Value = 33184.332
$ws.Range("K49").Value = 32249.75
$ws.Range("L49").Value = 33184.332
$ws.Range("M49").Value = -32019.75
$ws.Range("N49").Value = -33644.332
$ws.Range("H132").Value = 5820.7812
$ws.Range("I132").Value = 4822.2
$ws.Range("K132").Value = 14466.6
$ws.Range("M132").Value = -11936.6
